# Update countries & provincias Spain
# Applies the "paises.xlsx" data refresh: two countries swap rank position
# (their names swap rows while numbers are refreshed), the "datos
# actualizados" timestamp moves from 11:59 to 13:16, and a batch of
# countries get refreshed case/recovery/death counts for the newer
# snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header timestamp (row 1) -> updated snapshot time
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Octubre de 2020 a las 13:16"

# ---------------------------------------------------------------------
# Bonaire, San Eustaquio y Saba overtakes Liechtenstein in rank -> the
# two country names swap rows (195/196); numbers are then refreshed to
# the newer snapshot values.
# ---------------------------------------------------------------------
$ws.Range("A195").Value = "Bonaire, San Eustaquio y Saba"
$ws.Range("A196").Value = "Liechtenstein"

$ws.Range("B195").Value = 141
$ws.Range("C195").Value = 17
$ws.Range("D195").Value = 67
$ws.Range("E195").Value = 72
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 2

$ws.Range("B196").Value = 130
$ws.Range("C196").Value = 3
$ws.Range("D196").Value = 116
$ws.Range("E196").Value = 13
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 1

# ---------------------------------------------------------------------
# Islas Malvinas overtakes Montserrat in rank -> the two country names
# swap rows (215/216) along with their figures.
# ---------------------------------------------------------------------
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

$ws.Range("B216").Value = 13
$ws.Range("C216").Value = 0
$ws.Range("D216").Value = 12
$ws.Range("E216").Value = 0
$ws.Range("F216").Value = 0
$ws.Range("G216").Value = 0
$ws.Range("H216").Value = 1

# ---------------------------------------------------------------------
# Refreshed counts (Casos totales / Nuevos casos / Casos activos /
# Recuperados / Casos criticos / Muertes hoy / Muertes) for the rest of
# the updated countries, columns B:H.
# ---------------------------------------------------------------------

# India (row 5)
$ws.Range("B5").Value = 6687247
$ws.Range("C5").Value = 5174
$ws.Range("E5").Value = 921128
$ws.Range("G5").Value = 29
$ws.Range("H5").Value = 103629

# Iran (row 16)
$ws.Range("B16").Value = 479825
$ws.Range("C16").Value = 4151
$ws.Range("D16").Value = 394800
$ws.Range("E16").Value = 57606
$ws.Range("G16").Value = 227
$ws.Range("H16").Value = 27419

# Alemania (row 26)
$ws.Range("B26").Value = 304747
$ws.Range("C26").Value = 90
$ws.Range("E26").Value = 31430
$ws.Range("G26").Value = 1
$ws.Range("H26").Value = 9617

# Rumania (row 32)
$ws.Range("B32").Value = 139612
$ws.Range("C32").Value = 2121
$ws.Range("D32").Value = 109898
$ws.Range("E32").Value = 24593
$ws.Range("G32").Value = 73
$ws.Range("H32").Value = 5121

# Emiratos Arabes Unidos (row 44)
$ws.Range("B44").Value = 100794
$ws.Range("C44").Value = 1061
$ws.Range("D44").Value = 90556
$ws.Range("E44").Value = 9803
$ws.Range("G44").Value = 6
$ws.Range("H44").Value = 435

# Nepal (row 47)
$ws.Range("B47").Value = 90814
$ws.Range("C47").Value = 1551
$ws.Range("D47").Value = 67542
$ws.Range("E47").Value = 22709
$ws.Range("G47").Value = 9
$ws.Range("H47").Value = 563

# Suiza (row 62)
$ws.Range("B62").Value = 56632
$ws.Range("C62").Value = 700
$ws.Range("E62").Value = 7253
$ws.Range("G62").Value = 1
$ws.Range("H62").Value = 2079

# Senegal (row 93)
$ws.Range("B93").Value = 15141
$ws.Range("C93").Value = 19
$ws.Range("D93").Value = 12936
$ws.Range("E93").Value = 1893

# Uganda (row 110)
$ws.Range("B110").Value = 9082
$ws.Range("C110").Value = 117
$ws.Range("D110").Value = 5457
$ws.Range("E110").Value = 3541
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 84

# Malta (row 144)
$ws.Range("B144").Value = 3374
$ws.Range("C144").Value = 47
$ws.Range("D144").Value = 2812
$ws.Range("E144").Value = 522
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 40

# Gibraltar (row 181)
$ws.Range("B181").Value = 437
$ws.Range("C181").Value = 5
$ws.Range("D181").Value = 368
